$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: D3 changes from text "65674537" to a numeric value 65674537
$ws.Range("D3").Value = 65674537

# Row 4: new data ("SJ4" / "vijay 2" / "jhgyg@jjhj" / 1214512)
$ws.Range("A4").Value = "SJ4"
$ws.Range("B4").Value = "vijay 2"
$ws.Range("C4").Value = "jhgyg@jjhj"
$ws.Range("D4").Value = 1214512

# Row 5: new data ("SJ5" / "rich" / "r16@gmail.com" / 879545756454)
$ws.Range("A5").Value = "SJ5"
$ws.Range("B5").Value = "rich"
$ws.Range("C5").Value = "r16@gmail.com"
$ws.Range("D5").Value = 879545756454

# Row 6: new data ("SJ6" / "thullu" / "dsfd@dfs.com" / "45468645454")
# D6 must stay text (not be auto-converted to a number), so we briefly force
# a text format, assign the value, then restore the cell style so no lasting
# number-format change remains applied to the cell.
$ws.Range("A6").Value = "SJ6"
$ws.Range("B6").Value = "thullu"
$ws.Range("C6").Value = "dsfd@dfs.com"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "45468645454"
$ws.Range("D6").Style = "Normal"
